# Updated training results: replace ensemble metrics with new LM run.
# Every model row now carries the same (new) metric values, and the
# row labels in column A are re-shuffled per the new training order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New labels for column A, rows 2-26 (row 14 keeps "model_38_4_12")
$labels = @{
    2  = "model_38_4_0"
    3  = "model_38_4_22"
    4  = "model_38_4_21"
    5  = "model_38_4_20"
    6  = "model_38_4_19"
    7  = "model_38_4_18"
    8  = "model_38_4_17"
    9  = "model_38_4_16"
    10 = "model_38_4_15"
    11 = "model_38_4_14"
    12 = "model_38_4_13"
    13 = "model_38_4_23"
    14 = "model_38_4_12"
    15 = "model_38_4_10"
    16 = "model_38_4_9"
    17 = "model_38_4_8"
    18 = "model_38_4_7"
    19 = "model_38_4_6"
    20 = "model_38_4_5"
    21 = "model_38_4_4"
    22 = "model_38_4_3"
    23 = "model_38_4_2"
    24 = "model_38_4_1"
    25 = "model_38_4_11"
    26 = "model_38_4_24"
}

# New metric values (B..Q) shared by every data row (2-26)
# (values given as strings, then cast to double, so scientific
# notation literals parse correctly)
$values = @(
    [double]"0.9999106709542923",
    [double]"0.9989035467838921",
    [double]"0.9998688898941104",
    [double]"0.9998124870067278",
    [double]"0.9998620467351937",
    [double]"8.338474307615185e-05",
    [double]"0.001023490948502973",
    [double]"0.0001571851392513486",
    [double]"8.304504494770721e-05",
    [double]"0.0001201147639809278",
    [double]"0.0005467859721176896",
    [double]"0.009131524685185483",
    [double]"1.000064966578696",
    [double]"0.009520272595895827",
    [double]"132.7840904041427",
    [double]"202.2600124216301"
)

for ($row = 2; $row -le 26; $row++) {
    $ws.Cells.Item($row, 1).Value = $labels[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $values[$i]
    }
}
